# Apply the changes described by the diff:
#  - datetimeFigureOut placeholder text "9/6/14" -> "1/23/14" on the slide
#    master and on every slide layout (12 occurrences total).
#  - "HDFS/S3" -> "HDFS" and "Kinesis" -> "ZeroMQ" on slide 1's
#    architecture-diagram shapes (nested inside a group).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false

        if ($sh.Type -eq 14) {
            # msoPlaceholder
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    # ppPlaceholderDate
                    $isDatePlaceholder = $true
                }
            } catch {
            }
        }

        if (-not $isDatePlaceholder -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "9/6/14") {
                $isDatePlaceholder = $true
            }
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "9/6/14") {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes "1/23/14"

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "1/23/14"
}

# Slide 1: the streaming-sources diagram is a deeply nested group; PowerPoint
# flattens GroupItems so the member shapes are reachable directly.
$slide1 = $p.Slides.Item(1)
$topGroup = $slide1.Shapes.Item(1)
$items = $topGroup.GroupItems
for ($i = 1; $i -le $items.Count; $i++) {
    $sh = $items.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "HDFS/S3") {
            $tr.Text = "HDFS"
        } elseif ($tr.Text -eq "Kinesis") {
            $tr.Text = "ZeroMQ"
        }
    }
}
